$d = $word.ActiveDocument

# Locate the "Data, Technology and Strategy Consulting" paragraph under the
# Siege Analytics / PARTNER entry so we can insert the new bullet points
# right after it (and before the existing "Lead comprehensive polling..." bullet).
$idx = 0
$anchorIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t -match "^Data, Technology and Strategy Consulting") {
        $anchorIdx = $idx
        break
    }
}

if ($anchorIdx -eq -1) {
    throw "Could not find 'Data, Technology and Strategy Consulting' paragraph"
}

$newBullets = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

foreach ($line in $newBullets) {
    $p = $d.Paragraphs.Item($anchorIdx)
    $p.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($anchorIdx + 1)
    $newP.Range.Text = $line
    $anchorIdx = $anchorIdx + 1
}
